$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.0005009091434577991
$ws.Range("E2").Value = 0.0005009091434577991
$ws.Range("D3").Value = 0.5586335071389245
$ws.Range("E3").Value = 0.5586335071389245
$ws.Range("D4").Value = 0.001226009107456918
$ws.Range("E4").Value = 0.001226009107456918
$ws.Range("D5").Value = [double]"1.18235264240131E-10"
$ws.Range("E5").Value = [double]"1.18235264240131E-10"
$ws.Range("D6").Value = 0.801903460880978
$ws.Range("E6").Value = 0.801903460880978
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("C8").Value = $False
$ws.Range("D8").Value = [double]"1.789834212173071E-07"
$ws.Range("E8").Value = 0.9999998210165788
$ws.Range("C9").Value = $False
$ws.Range("D9").Value = 0.09986722506612214
$ws.Range("E9").Value = 0.9001327749338779
$ws.Range("C10").Value = $False
$ws.Range("D10").Value = [double]"2.169052331615286E-05"
$ws.Range("E10").Value = 0.9999783094766839
$ws.Range("D11").Value = 0.9999572275875143
$ws.Range("E11").Value = [double]"4.277241248573649E-05"
$ws.Range("F11").Value = 3.101717472076416
$ws.Range("G11").Value = 0.5
$ws.Range("D12").Value = [double]"7.733603271807684E-06"
$ws.Range("E12").Value = [double]"7.733603271807684E-06"
$ws.Range("D13").Value = 0.7052595744657145
$ws.Range("E13").Value = 0.7052595744657145
$ws.Range("D14").Value = [double]"1.436017113963682E-06"
$ws.Range("E14").Value = [double]"1.436017113963682E-06"
$ws.Range("D15").Value = [double]"8.843900312332742E-15"
$ws.Range("E15").Value = [double]"8.843900312332742E-15"
$ws.Range("D16").Value = 0.6727888486231047
$ws.Range("E16").Value = 0.6727888486231047
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("C18").Value = $False
$ws.Range("D18").Value = [double]"2.006023524455909E-11"
$ws.Range("E18").Value = 0.9999999999799397
$ws.Range("C19").Value = $False
$ws.Range("D19").Value = 0.05143777588323453
$ws.Range("E19").Value = 0.9485622241167655
$ws.Range("C20").Value = $False
$ws.Range("D20").Value = [double]"1.006920530587455E-08"
$ws.Range("E20").Value = 0.9999999899307946
$ws.Range("D21").Value = 0.9999832914291426
$ws.Range("E21").Value = [double]"1.670857085744881E-05"
$ws.Range("F21").Value = 4.83522891998291
$ws.Range("G21").Value = 0.5
